$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.901.47"
$ws.Cells.Item(2, 5).Value = "  +0.78%  "

$ws.Cells.Item(3, 4).Value = "2.321.25"
$ws.Cells.Item(3, 5).Value = "  +1.63%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "302.41"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -0.46%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "96.36"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +0.75%  "

$ws.Cells.Item(7, 5).Value = "  +0.72%  "

$ws.Cells.Item(8, 5).Value = "  +0.06%  "

$ws.Cells.Item(9, 5).Value = "  +0.41%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "34.63"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -0.25%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "19.01"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +5.85%  "

$ws.Cells.Item(12, 5).Value = "  +0.65%  "

$ws.Cells.Item(13, 5).Value = "  +0.44%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.78"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +0.21%  "

$ws.Cells.Item(15, 4).Value = "2.686.21"
$ws.Cells.Item(15, 5).Value = "  +1.77%  "

$ws.Cells.Item(16, 4).Value = "2.322.40"
$ws.Cells.Item(16, 5).Value = "  +1.32%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.791"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +2.53%  "

$ws.Cells.Item(18, 4).Value = "42.832.44"
$ws.Cells.Item(18, 5).Value = "  +0.84%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.20"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -5.67%  "

$ws.Cells.Item(20, 5).Value = "  +3.53%  "

$ws.Cells.Item(21, 5).Value = "  +0.36%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "68.00"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +1.35%  "

$ws.Cells.Item(23, 5).Value = "  +6.88%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "236.38"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +0.26%  "

$ws.Cells.Item(25, 5).Value = "  -0.03%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.43"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +1.21%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.45"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -0.82%  "

$ws.Cells.Item(28, 5).Value = "  -1.04%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "166.38"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +0.07%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.14"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +1.99%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "32.37"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -0.56%  "

$ws.Cells.Item(32, 5).Value = "  -0.01%  "

$ws.Cells.Item(33, 5).Value = "  +1.50%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "17.84"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +0.24%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.48"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +1.41%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.32"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -1.39%  "

$ws.Cells.Item(38, 5).Value = "  +3.77%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0998"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -0.24%  "

$ws.Cells.Item(40, 5).Value = "  +3.75%  "

$ws.Cells.Item(41, 5).Value = "  +0.46%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "21.02"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +14.45%  "

$ws.Cells.Item(43, 4).Value = "1.936.86"
$ws.Cells.Item(43, 5).Value = "  -2.65%  "

$ws.Cells.Item(44, 5).Value = "  +1.47%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.23"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +0.61%  "

$ws.Cells.Item(46, 5).Value = "  +4.37%  "

$ws.Cells.Item(47, 5).Value = "  +1.03%  "

$ws.Cells.Item(48, 4).Value = "2.553.27"
$ws.Cells.Item(48, 5).Value = "  +1.88%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "53.52"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -0.05%  "

$ws.Cells.Item(50, 5).Value = "  -3.58%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "72.14"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +2.55%  "
